$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "288.54"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-9.81%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.49%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.035"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.37%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07297"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.78%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.285"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.28%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.520"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-14.06%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9192"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.60%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1187"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.46%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1731"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-7.86%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08620"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04168"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.25%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.15%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001266"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.32%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005842"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.45%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.72%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3290"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.04%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.864"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.81%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.33%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03863"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.15%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.08%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003814"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-7.43%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001281"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.85%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003725"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02329"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-8.86%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.04973"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.99%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006346"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "218.80%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007683"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.06%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1274"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.37%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007346"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.20%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007067"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.60%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3128"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.52%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006439"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.90%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08746"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-56.14%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.13%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.09%"
